$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15
$ws.Cells.Item($row, 1).Value = "2025-09-23 05:00:13"
$ws.Cells.Item($row, 2).Value = "Gradle Basics"
$ws.Cells.Item($row, 3).Value = "Build Tools"
$ws.Cells.Item($row, 4).Value = 486
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = "SUCCESS"
$ws.Cells.Item($row, 7).Value = "Generated successfully"
